$d = $word.ActiveDocument

$table = $d.Tables.Item(1)
$row = $table.Rows.Item(6)

$row.Cells.Item(1).Range.Text = "6/08/14"
$row.Cells.Item(2).Range.Text = "14:00-18:00"
$row.Cells.Item(3).Range.Text = "Developed fixes for scaling and GAP requests"
